$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.33), maa://25390 (95.98), maa://36681 (87.01)'
$ws.Range('P5').Value = 'maa://21919 (96.15), maa://21281 (85.71)'
$ws.Range('P6').Value = 'maa://31836 (91.67), maa://30381 (92.31)'
$ws.Range('A8').Value = '更新日期：2025.01.12 13:18:28'
$ws.Range('H8').Value = '*maa://24371 (54.93)'
$ws.Range('P8').Value = 'maa://32931 (84.11), *maa://21916 (61.29), maa://23252 (92.42), maa://37496 (96.55), **maa://22759 (45.45)'
$ws.Range('AF8').Value = '*maa://24479 (77.38), *maa://21990 (51.85)'
$ws.Range('D9').Value = 'maa://22765 (92.31), *maa://21915 (68.0)'
$ws.Range('P9').Value = 'maa://22736 (81.72)'
$ws.Range('T9').Value = '**maa://22866 (30.19), maa://26222 (97.83)'
$ws.Range('D10').Value = '***maa://25695 (19.34), **maa://32237 (41.86), ***maa://34206 (20.83), ***maa://39951 (16.28), ***maa://39243 (28.57), **maa://45271 (50.0)'
$ws.Range('P10').Value = 'maa://28977 (91.36), maa://36669 (89.74), *maa://23264 (61.82)'
$ws.Range('X10').Value = 'maa://22301 (97.67), maa://22726 (100.0)'
$ws.Range('L11').Value = 'maa://21287 (88.54)'
$ws.Range('D13').Value = 'maa://24999 (91.83), maa://36673 (92.86), maa://25001 (85.51)'
$ws.Range('P13').Value = 'maa://22676 (92.11), *maa://22583 (74.24), *maa://22500 (57.78)'
$ws.Range('D15').Value = '*maa://22743 (77.66), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (35.42), maa://45058 (100.0)'
$ws.Range('AF15').Value = 'maa://21364 (81.15), *maa://22766 (70.27), *maa://36666 (78.16)'
$ws.Range('D16').Value = 'maa://21441 (96.35), maa://36679 (93.48), maa://37650 (96.97)'
$ws.Range('H17').Value = 'maa://22430 (88.83), maa://39599 (87.18)'
$ws.Range('AF18').Value = '*maa://24313 (58.13), **maa://29784 (44.44)'
$ws.Range('D20').Value = 'maa://21432 (90.07), maa://25198 (93.14), *maa://20795 (51.56), maa://36680 (93.55)'
$ws.Range('L20').Value = 'maa://41331 (85.59)'
$ws.Range('D21').Value = 'maa://21261 (97.44)'
$ws.Range('L22').Value = 'maa://27127 (84.26), *maa://22751 (73.85)'
$ws.Range('D23').Value = '***maa://28036 (28.57), *maa://41753 (53.85)'
$ws.Range('L23').Value = 'maa://39756 (94.58), maa://39875 (93.85)'
$ws.Range('D24').Value = '*maa://24368 (78.36)'
$ws.Range('D25').Value = 'maa://29753 (94.96)'
$ws.Range('H25').Value = '*maa://29063 (74.17), *maa://25311 (73.53), ***maa://22725 (4.84), maa://45047 (100.0)'
$ws.Range('AB25').Value = 'maa://31215 (86.27), *maa://24516 (79.78), maa://26001 (87.5)'
$ws.Range('AF25').Value = 'maa://20108 (96.27), maa://24621 (96.61), maa://36676 (96.88), maa://22771 (85.71), **maa://37772 (50.0)'
$ws.Range('D28').Value = 'maa://24465 (90.81), maa://25725 (83.72)'
$ws.Range('AF28').Value = 'maa://36660 (93.06), *maa://36701 (64.29)'
$ws.Range('AB30').Value = 'maa://42979 (96.69), maa://45045 (100.0)'
$ws.Range('H32').Value = 'maa://21895 (97.41), maa://36667 (98.44), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (95.6), maa://41108 (88.0), maa://41238 (96.43)'
$ws.Range('L35').Value = 'maa://41296 (96.85)'
$ws.Range('AF35').Value = 'maa://39479 (93.75)'
$ws.Range('P37').Value = 'maa://21280 (88.83), *maa://21239 (66.67)'
$ws.Range('H39').Value = 'maa://25199 (84.82), maa://36670 (87.91), maa://30434 (89.86), ***maa://25036 (16.0), *maa://44165 (66.67), *maa://45059 (80.0)'
$ws.Range('P39').Value = 'maa://24709 (91.67)'
$ws.Range('P40').Value = 'maa://23278 (95.68), maa://21386 (95.74), maa://36664 (90.91)'
$ws.Range('H45').Value = 'maa://21229 (84.86), maa://30807 (95.52), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (81.82)'
$ws.Range('H46').Value = 'maa://35931 (92.41), maa://43901 (91.67)'
$ws.Range('H55').Value = 'maa://32532 (92.08)'
